$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.544.26'
$ws.Range("E2").Value = '  +2.35%  '
$ws.Range("D3").Value = '1.684.48'
$ws.Range("E3").Value = '  +3.25%  '
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '217.53'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +5.08%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.5335'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +3.90%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -0.10%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.2685'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +5.36%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.06435'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +4.13%  '
$ws.Range("E10").Value = '  +6.35%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.07788'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +3.53%  '
$ws.Range("D12").Value = '1.685.86'
$ws.Range("E12").Value = '  +3.10%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '4.503'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +3.78%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.5624'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +5.15%  '
$ws.Range("D15").Value = '0.0₅8439'
$ws.Range("E15").Value = '  +6.63%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '66.07'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +2.17%  '
$ws.Range("D17").Value = '26.588.36'
$ws.Range("E17").Value = '  +2.53%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -0.15%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '4.817'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +4.57%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '194.54'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +5.41%  '
$ws.Range("E21").Value = '  +5.26%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '6.400'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +6.06%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '1.003'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -0.07%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '144.02'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -1.31%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.1272'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +6.98%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '7.483'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +2.95%  '
$ws.Range("E27").Value = '  +5.41%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '1.416'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +3.37%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '0.06122'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +3.22%  '
$ws.Range("E30").Value = '  +3.10%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '3.608'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +8.04%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '3.464'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +4.30%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.700'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +6.93%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.019'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +5.98%  '
$ws.Range("E35").Value = '  +2.85%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '2.419'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +1.74%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.5714'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.92%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.01641'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +4.16%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '5.979'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +4.06%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.8683'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +3.96%  '
$ws.Range("D41").Value = '1.057.18'
$ws.Range("E41").Value = '  +0.65%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '100.13'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").Value = '1.835.47'
$ws.Range("E44").Value = '  +2.91%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '57.23'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +6.24%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '8.159'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +3.68%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -0.24%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.05209'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +0.32%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '6.064'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +5.41%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.4243'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +0.49%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.09930'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +4.72%  '
